$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Serial No" record row (currently row 64, A64:G64) needs to move down
# past the "App Version No" (65) and "Battery Level" (66) rows, ending up
# in row 67 (the row right before "Bluetooth Name", which stays at row 68).
#
# Net effect:
#   before: 64=Serial No, 65=App Version No, 66=Battery Level, 67=<blank>, 68=Bluetooth Name
#   after:  64=App Version No, 65=Battery Level, 66=<blank>, 67=Serial No, 68=Bluetooth Name
#
# Implemented as: insert a fresh blank row just above "Bluetooth Name" (68),
# copy the "Serial No" row's values+formatting into that new blank row, then
# delete the original "Serial No" row (64), which shifts "App Version No" and
# "Battery Level" up by one and leaves the blank separator row in its new spot.

$ws.Rows.Item(68).Insert()
$ws.Range("A64:G64").Copy($ws.Range("A68:G68"))
$ws.Rows.Item(64).Delete()

# Update the sheet view to match where the user ended up looking after the edit.
$aw = $excel.ActiveWindow
$ws.Range("G67").Select()
$aw.ScrollRow = 31
$aw.ScrollColumn = 1
